# Add a new "chucVu" (role) column as column B, shifting the existing
# columns (taiKhoan..email) one position to the right, and populate it
# with "admin" for the two admin accounts (rows 2 and 38) and "user"
# for everyone else.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before the existing column B ("taiKhoan"),
# pushing every other column (B..H) one slot to the right (C..I).
$ws.Columns("B:B").Insert()

# Header for the new column.
$ws.Range("B1").Value2 = "chucVu"

# Row 2 (maNVYT 1605) is an admin account.
$ws.Cells.Item(2, 2).Value2 = "admin"

# Default every other data row (3..54) to "user" first so the shared
# string table picks up "user" before we revisit row 38 below.
for ($r = 3; $r -le 54; $r++) {
    $ws.Cells.Item($r, 2).Value2 = "user"
}

# Row 38 (maNVYT 1173) is the other admin account.
$ws.Cells.Item(38, 2).Value2 = "admin"

# Restore the cursor/selection position recorded in the saved file.
$ws.Range("C17").Select()
